$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.042494666666667
$ws.Range("H2").Value = 3.127484
$ws.Range("I2").Value = 0.0007670466909205676
$ws.Range("J2").Value = 0.0007670466909205677
$ws.Range("M2").Value = 83.91225566666667
$ws.Range("N2").Value = 251.736767
$ws.Range("O2").Value = 0.9556261553553385
$ws.Range("P2").Value = 0.9556261553553385
$ws.Range("Q2").Value = 87.47807900046978
$ws.Range("R2").Value = 787.3027110042281
$ws.Range("S2").Value = 0.0007330098802224566
$ws.Range("T2").Value = 0.0007330098802224567

# Row 3
$ws.Range("G3").Value = 1.042494666666667
$ws.Range("H3").Value = 3.127484
$ws.Range("I3").Value = 0.0007670466909205676
$ws.Range("J3").Value = 0.0007670466909205677
$ws.Range("O3").Value = 0.00439999103960854
$ws.Range("P3").Value = 0.00439999103960854
$ws.Range("Q3").Value = 0.4027754594275555
$ws.Range("R3").Value = 3.624979134848
$ws.Range("S3").Value = 0.000003374998567011879
$ws.Range("T3").Value = 0.000003374998567011879

# Row 4
$ws.Range("G4").Value = 1.042494666666667
$ws.Range("H4").Value = 3.127484
$ws.Range("I4").Value = 0.0007670466909205676
$ws.Range("J4").Value = 0.0007670466909205677
$ws.Range("M4").Value = 3.510050666666667
$ws.Range("N4").Value = 10.530152
$ws.Range("O4").Value = 0.03997385360505296
$ws.Range("P4").Value = 0.03997385360505297
$ws.Range("Q4").Value = 3.659209099729777
$ws.Range("R4").Value = 32.932881897568
$ws.Range("S4").Value = 0.00003066181213109908
$ws.Range("T4").Value = 0.00003066181213109909

# Row 5
$ws.Range("I5").Value = 0.9658609009611662
$ws.Range("J5").Value = 0.9658609009611662
$ws.Range("M5").Value = 83.91225566666667
$ws.Range("N5").Value = 251.736767
$ws.Range("O5").Value = 0.9556261553553385
$ws.Range("P5").Value = 0.9556261553553385
$ws.Range("Q5").Value = 110151.9075668569
$ws.Range("R5").Value = 991367.168101712
$ws.Range("S5").Value = 0.9230019393935626
$ws.Range("T5").Value = 0.9230019393935626

# Row 6
$ws.Range("I6").Value = 0.9658609009611662
$ws.Range("J6").Value = 0.9658609009611662
$ws.Range("O6").Value = 0.00439999103960854
$ws.Range("P6").Value = 0.00439999103960854
$ws.Range("S6").Value = 0.004249779309737363
$ws.Range("T6").Value = 0.004249779309737363

# Row 7
$ws.Range("I7").Value = 0.9658609009611662
$ws.Range("J7").Value = 0.9658609009611662
$ws.Range("M7").Value = 3.510050666666667
$ws.Range("N7").Value = 10.530152
$ws.Range("O7").Value = 0.03997385360505296
$ws.Range("P7").Value = 0.03997385360505297
$ws.Range("Q7").Value = 4607.655622148167
$ws.Range("R7").Value = 41468.90059933351
$ws.Range("S7").Value = 0.03860918225786621
$ws.Range("T7").Value = 0.03860918225786622

# Row 8
$ws.Range("G8").Value = 45.356022
$ws.Range("H8").Value = 136.068066
$ws.Range("I8").Value = 0.03337205234791334
$ws.Range("J8").Value = 0.03337205234791334
$ws.Range("M8").Value = 83.91225566666667
$ws.Range("N8").Value = 251.736767
$ws.Range("O8").Value = 0.9556261553553385
$ws.Range("P8").Value = 0.9556261553553385
$ws.Range("Q8").Value = 3805.926114086958
$ws.Range("R8").Value = 34253.33502678262
$ws.Range("S8").Value = 0.03189120608155352
$ws.Range("T8").Value = 0.03189120608155352

# Row 9
$ws.Range("G9").Value = 45.356022
$ws.Range("H9").Value = 136.068066
$ws.Range("I9").Value = 0.03337205234791334
$ws.Range("J9").Value = 0.03337205234791334
$ws.Range("O9").Value = 0.00439999103960854
$ws.Range("P9").Value = 0.00439999103960854
$ws.Range("Q9").Value = 17.523631710528
$ws.Range("R9").Value = 157.712685394752
$ws.Range("S9").Value = 0.0001468367313041658
$ws.Range("T9").Value = 0.0001468367313041658

# Row 10
$ws.Range("G10").Value = 45.356022
$ws.Range("H10").Value = 136.068066
$ws.Range("I10").Value = 0.03337205234791334
$ws.Range("J10").Value = 0.03337205234791334
$ws.Range("M10").Value = 3.510050666666667
$ws.Range("N10").Value = 10.530152
$ws.Range("O10").Value = 0.03997385360505296
$ws.Range("P10").Value = 0.03997385360505297
$ws.Range("Q10").Value = 159.201935258448
$ws.Range("R10").Value = 1432.817417326032
$ws.Range("S10").Value = 0.001334009535055652
$ws.Range("T10").Value = 0.001334009535055652
